# Add "2022-Q1" fund-holdings sheet (positioned between "2021-Q4" and
# "总计") and update the "总计" (summary) sheet with a new row for the
# 2022-Q1 quarter.
#
# Strategy: the existing "总计" worksheet is renamed to "2022-Q1" and
# re-populated with the fund-holdings table (this keeps the original
# sheetId/relationship slot with the new quarterly data, matching how
# the quarterly sheets are normally appended). A brand-new worksheet is
# then added right after it and named "总计", holding the refreshed
# summary table (2022-Q1 plus the four previously existing quarters).

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$oldTotal.Name = "2022-Q1"
$newQ1 = $oldTotal
$newQ1.Cells.Clear()

# Copy header (bold/border/center) formatting from the 2021-Q4 sheet so
# the new sheet matches the look of the other quarterly sheets.
$q4.Range("B1:H1").Copy()
$newQ1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$newQ1.Range("A2:A16").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $newQ1.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# code, name, scale, stockPosition, positionRatio, marketValue, mvIsNumber, rank
$fundData = @(
    ,@("470009","汇添富民营活力混合A","25.42","88.81","5.22","1.3269",0,5)
    ,@("007355","汇添富科技创新灵活配置混合A","16.80","85.87","4.44","0.7459",0,5)
    ,@("690011","民生加银积极成长混合","4.04","85.33","6.39","0.2582",0,5)
    ,@("000884","民生加银优选股票","2.09","84.80","6.38","0.1333",0,5)
    ,@("007356","汇添富科技创新灵活配置混合C","2.52","85.87","4.44","0.1119",0,5)
    ,@("013296","民生加银聚优精选混合","1.41","84.53","6.28","0.0885",0,4)
    ,@("011888","民生加银周期优选混合型证券投资基金A","0.54","86.02","6.52","0.0352",0,4)
    ,@("004194","招商中证1000指数增强A","1.76","94.40","1.10","0.0194",0,3)
    ,@("011889","民生加银周期优选混合型证券投资基金C","0.12","86.02","6.52","0.0078",0,4)
    ,@("004195","招商中证1000指数增强C","0.68","94.40","1.10","0.0075",0,3)
    ,@("004726","先锋聚优灵活配置混合A","0.06","93.23","5.08","0.0030",0,7)
    ,@("004727","先锋聚优灵活配置混合C","0.04","93.23","5.08","0.0020",0,7)
    ,@("004833","先锋聚利灵活配置混合A","0.02","94.68","4.56","0.0009",0,10)
    ,@("004834","先锋聚利灵活配置混合C","0.02","94.68","4.56","0.0009",0,10)
    ,@("960014","汇添富民营活力混合型证券投资基金 O","0.00","88.81","5.22",0,1,5)
)

# Text-like columns (B..G) need to be forced to Text format first,
# otherwise Excel auto-converts numeric-looking strings (fund codes,
# percentages, market values) to numbers and e.g. loses leading zeros.
$newQ1.Range("B2:G16").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $row = 2 + $i
    $rec = $fundData[$i]
    $newQ1.Cells.Item($row, 1).Value = $i
    $newQ1.Cells.Item($row, 2).Value = $rec[0]
    $newQ1.Cells.Item($row, 3).Value = $rec[1]
    $newQ1.Cells.Item($row, 4).Value = $rec[2]
    $newQ1.Cells.Item($row, 5).Value = $rec[3]
    $newQ1.Cells.Item($row, 6).Value = $rec[4]
    if ($rec[6] -eq 1) {
        $mvCell = $newQ1.Cells.Item($row, 7)
        $mvCell.NumberFormat = "General"
        $mvCell.Value = [double]$rec[5]
    } else {
        $newQ1.Cells.Item($row, 7).Value = $rec[5]
    }
    $newQ1.Cells.Item($row, 8).Value = $rec[7]
}

# ---------------------------------------------------------------------
# Step 2: create a brand-new "总计" sheet right after "2022-Q1" holding
# the refreshed summary table.
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $newQ1)
$newTotal.Name = "总计"

$newQ1.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$newQ1.Range("A2").Copy()
$newTotal.Range("A2:A6").PasteSpecial(-4122)

$newTotal.Cells.Item(1, 2).Value = "日期"
$newTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$newTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summary = @(
    ,@("2022-Q1", 15, 2.74)
    ,@("2021-Q4", 10, 4.74)
    ,@("2021-Q3", 2, 1.25)
    ,@("2021-Q2", 8, 2.78)
    ,@("2021-Q1", 2, 0.05)
)

for ($i = 0; $i -lt $summary.Count; $i++) {
    $row = 2 + $i
    $rec = $summary[$i]
    $newTotal.Cells.Item($row, 1).Value = $i
    $newTotal.Cells.Item($row, 2).Value = $rec[0]
    $newTotal.Cells.Item($row, 3).Value = $rec[1]
    $newTotal.Cells.Item($row, 4).Value = $rec[2]
}

# Restore the original active sheet/selection (adding worksheets makes
# the newest one active; the source workbook had "2021-Q1" selected).
$wb.Worksheets.Item("2021-Q1").Activate()
